## Update on 20181114.2334 by YKBKyle
## Rebuild the "RNote" worksheet: rename sheet, rewrite the R-command
## reference table (new rows for merge/select/filter/packageVersion),
## add a hyperlink + merged label cells, resize columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Rename the worksheet
# ---------------------------------------------------------------
$ws.Name = "RCommands"

# ---------------------------------------------------------------
# 2. Cell values (row by row, left to right)
# ---------------------------------------------------------------
$ws.Range("A1").Value2 = "Commands"
$ws.Range("B1").Value2 = "Example"
$ws.Range("C1").Value2 = "Description"
$ws.Range("D1").Value2 = "Learning Resources"
$ws.Range("E1").Value2 = "Class"

$ws.Range("A2").Value2 = "rm(list=ls())"
$ws.Range("C2").Value2 = "remove all the data in the global environment"
$ws.Range("E2").Value2 = "Other"

$ws.Range("A3").Value2 = 'install.packages("Package")'
$ws.Range("C3").Value2 = "install package"
$ws.Range("E3").Value2 = "Package Operation"
# D3 used to hold "package operation" - no longer used, drop it entirely
$ws.Range("D3").Clear()

$ws.Range("A4").Value2 = 'packageVersion("Package")'
$ws.Range("C4").Value2 = "check the version of the loaded package"
$ws.Range("E4").Value2 = "Package Operation"

$ws.Range("A5").Value2 = "merge()"
$ws.Range("B5").Value2 = 'total <- merge(data frame A, data frame B, by = "ID")'
$ws.Range("C5").Value2 = "merge 2 data frames by ID"
$ws.Range("D5").Value2 = "http://www.statmethods.net/management/merging.html"
$ws.Range("E5").Value2 = "Data Management"

$ws.Range("B6").Value2 = 'total <- merge(data frame A, data frame B, by =c("ID","Country"))'
$ws.Range("C6").Value2 = "merge 2 data frames by ID and Country"

$ws.Range("A7").Value2 = "tbl_df()"
$ws.Range("B7").Value2 = "tbl_df(data frame A)"
$ws.Range("C7").Value2 = "convert date frame to tibble(data frame tbl)"
$ws.Range("E7").Value2 = "Package: dplyr"

$ws.Range("A8").Value2 = "select()"
$ws.Range("B8").Value2 = "select(data frame A, ColVar1, ColVar2, ColVar3)"
$ws.Range("C8").Value2 = "subset only 3 columns of A"
$ws.Range("E8").Value2 = "Package: dplyr"

$ws.Range("B9").Value2 = "select(data frame A, ColVar1:ColVarN)"
$ws.Range("C9").Value2 = " subset all columns starting from ColVar1 and ending with ColVarN"

$ws.Range("B10").Value2 = "select(data frame A, -(ColVar1:ColVarN))"
$ws.Range("C10").Value2 = "subset all columns except those from ColVar1 to ColVarN"

$ws.Range("A11").Value2 = "filter()"
$ws.Range("B11").Value2 = "filter(data frame A, ColVar > 0)"
$ws.Range("C11").Value2 = "subset all rows where ColVar > 0"
$ws.Range("E11").Value2 = "Package: dplyr"

# ---------------------------------------------------------------
# 3. Formatting: the "group label" cells (A5,E5 / A8,E8 and the
#    blank continuation cells underneath them) use the same
#    centered 14pt Arial look as the header/example rows.
# ---------------------------------------------------------------
$labelCells = @("A5","E5","A6","E6","A8","E8","A9","E9","A10","E10")
foreach ($ref in $labelCells) {
    $rng = $ws.Range($ref)
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 14
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}

# ---------------------------------------------------------------
# 4. Merge the label cells
# ---------------------------------------------------------------
$ws.Range("A5:A6").Merge()
$ws.Range("E5:E6").Merge()
$ws.Range("A8:A10").Merge()
$ws.Range("E8:E10").Merge()

# ---------------------------------------------------------------
# 5. Hyperlink on D5 (merged with D6)
# ---------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D5"), "http://www.statmethods.net/management/merging.html")
$ws.Range("D5").HorizontalAlignment = -4108
$ws.Range("D5").VerticalAlignment = -4108
$ws.Range("D5:D6").Merge()

# ---------------------------------------------------------------
# 6. Column widths (character units; engine stores width in 1/7
#    increments plus the standard 5/7 padding, so feed it the
#    inverse so the saved XML width lands on the target value).
# ---------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 67.61774553571429
$ws.Columns.Item(3).ColumnWidth = 69.61774553571429
$ws.Columns.Item(4).ColumnWidth = 21.617745535714285
$ws.Columns.Item(5).ColumnWidth = 19.949776785714285

# ---------------------------------------------------------------
# 7. Selection / active cell
# ---------------------------------------------------------------
$ws.Range("B12").Select()

# ---------------------------------------------------------------
# 8. Window position (best effort - engine may not persist this)
# ---------------------------------------------------------------
try { $excel.ActiveWindow.Left = 1380 } catch {}
try { $excel.ActiveWindow.Top = 4580 } catch {}
